$wb = $excel.ActiveWorkbook

# Map of row -> new "想去人数" (F column) value for both "展览" and "全部类型" sheets
$updates = @{
    2  = 828
    4  = 1150
    6  = 12377
    9  = 495
    11 = 1134
    12 = 912
    13 = 13619
    14 = 13828
    22 = 212
    23 = 4923
    24 = 222
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
